# Add data for 2021-11-08: update "through 10-30" -> "through 10-31"
# labels and refresh October / Total rows with the latest counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name shown in the workbook's sheet list.
$ws.Name = "Through 2021-10-31"

# Update the "October (through 10-30)" label in column A, row 11.
$ws.Range("A11").Value = "October (through 10-31)"

# Update October row (row 11) values for years 2015-2021 (columns B-H).
$ws.Range("B11").Value = 32
$ws.Range("C11").Value = 57
$ws.Range("D11").Value = 83
$ws.Range("E11").Value = 67
$ws.Range("F11").Value = 60
$ws.Range("G11").Value = 156
$ws.Range("H11").Value = 196

# Update Total row (row 12) values for years 2015-2021 (columns B-H).
$ws.Range("B12").Value = 258
$ws.Range("C12").Value = 486
$ws.Range("D12").Value = 710
$ws.Range("E12").Value = 615
$ws.Range("F12").Value = 482
$ws.Range("G12").Value = 1057
$ws.Range("H12").Value = 1444
